$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared / rich-text strings) ---
# "Volume 29   Number  44" -> "Volume 29   Number  45"
$ws.Range("A8").Value = "Volume 29   Number  45"
# "Report Covering the Week  10/31/2022  Through  11/6/2022"
#   -> "Report Covering the Week  11/7/2022  Through  11/13/2022"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Cells that change representation: numeric -> text "0" (shared string, style 14) ---
# Use a cell that already holds that exact text+style (D15) as the template:
#  1) PasteSpecial formats (xlPasteFormats = -4122) to adopt style 14
#  2) PasteSpecial all (xlPasteAll = -4104) to adopt the literal text value "0"
$fmtSrc = $ws.Range("D15")

$fmtSrc.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("C15").PasteSpecial(-4104)

$fmtSrc.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("C20").PasteSpecial(-4104)

$fmtSrc.Copy()
$ws.Range("C26").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("C26").PasteSpecial(-4104)

# --- Cells that change representation: text -> numeric (row 27) ---
# D27 was text "0" (style 14) -> becomes numeric 1 (style 15, like F27)
$ws.Range("F27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1

# E27 was text "***.*" (style 14) -> becomes numeric 0 (style 16, like H27)
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0

# --- Plain numeric value updates ---
# Row 15
$ws.Range("L15").Value = 5.882352941176
$ws.Range("M15").Value = -21.739130434782
$ws.Range("N15").Value = -68.965517241379

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 17
$ws.Range("H16").Value = -5.555555555555
$ws.Range("I16").Value = 171
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 6.875
$ws.Range("L16").Value = 20.422535211267
$ws.Range("M16").Value = -22.272727272727
$ws.Range("N16").Value = -78.678304239401

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -2.941176470588
$ws.Range("I17").Value = 337
$ws.Range("J17").Value = 322
$ws.Range("K17").Value = 4.658385093167
$ws.Range("L17").Value = 24.814814814814
$ws.Range("M17").Value = 74.611398963730
$ws.Range("N17").Value = -36.174242424242

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("I18").Value = 169
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 85.714285714285
$ws.Range("L18").Value = 18.181818181818
$ws.Range("M18").Value = 96.511627906976
$ws.Range("N18").Value = -58.880778588807

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 77.777777777777
$ws.Range("I19").Value = 429
$ws.Range("J19").Value = 317
$ws.Range("K19").Value = 35.331230283911
$ws.Range("L19").Value = 44.932432432432
$ws.Range("M19").Value = 88.986784140969
$ws.Range("N19").Value = -16.046966731898

# Row 20
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -28.571428571428
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = 39.622641509434
$ws.Range("L20").Value = 27.586206896551
$ws.Range("M20").Value = 117.647058823529
$ws.Range("N20").Value = -80.053908355795

# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 29.787234042553
$ws.Range("I21").Value = 1202
$ws.Range("J21").Value = 965
$ws.Range("K21").Value = 24.559585492228
$ws.Range("L21").Value = 28.693790149892
$ws.Range("M21").Value = 52.926208651399
$ws.Range("N21").Value = -55.596601403768

# Row 23
$ws.Range("C23").Value = 7
$ws.Range("E23").Value = -22.222222222222
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = -10.810810810810
$ws.Range("I23").Value = 360
$ws.Range("J23").Value = 372
$ws.Range("K23").Value = -3.225806451612
$ws.Range("L23").Value = -1.639344262295
$ws.Range("M23").Value = 45.748987854251

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 41.666666666666
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = 21.153846153846
$ws.Range("I24").Value = 775
$ws.Range("J24").Value = 666
$ws.Range("K24").Value = 16.366366366366
$ws.Range("L24").Value = 17.424242424242
$ws.Range("M24").Value = 32.027257240204

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = -40.983606557377
$ws.Range("I25").Value = 479
$ws.Range("J25").Value = 471
$ws.Range("K25").Value = 1.698513800424
$ws.Range("L25").Value = -0.415800415800
$ws.Range("M25").Value = -19.224283305227

# Row 26
$ws.Range("L26").Value = 8.333333333333

# Row 27 (remaining plain numeric updates)
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 58
$ws.Range("J27").Value = 66
$ws.Range("K27").Value = -12.121212121212
$ws.Range("L27").Value = 34.883720930232

# Row 28
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -20
$ws.Range("N28").Value = -73.626373626373

# Row 29
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 28
$ws.Range("K29").Value = -39.285714285714
$ws.Range("N29").Value = -80
